$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column B (shifts old B,C,D -> C,D,E)
$ws.Columns("B").Insert()

# 2. Header row (row 1)
$ws.Range("B1").Value = "contentStatus"

# 3. New "contentStatus" column values for data rows
$ws.Range("B2").Value = "standard"
$ws.Range("B3").Value = "standard"
$ws.Range("B4").Value = "standard"
$ws.Range("B5").Value = "standard"
$ws.Range("B6").Value = "standard"
$ws.Range("B7").Value = "standard"

# 4. New row 8 (duplicate of row 7, plain/no special formatting except a number format on E8)
$ws.Range("A8").Value = "Grupo de Datos 2"
$ws.Range("B8").Value = "standard"
$ws.Range("C8").Value = "Dato grupo 2"
$ws.Range("D8").Value = "Equipo veterano da un gran espectaculo"
$ws.Range("E8").Value = 45261
$ws.Range("E8").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 5. Column widths (target widths of 17.71/14.71/13.71/46.71/18.71 chars; the
# host engine snaps ColumnWidth to whole-pixel increments, so we feed the
# value that lands in the middle of the pixel bucket nearest the target)
$ws.Columns("A").ColumnWidth = 16.8333333
$ws.Columns("B").ColumnWidth = 13.8333333
$ws.Columns("C").ColumnWidth = 12.8333333
$ws.Columns("D").ColumnWidth = 45.8333333
$ws.Columns("E").ColumnWidth = 17.8333333

# 6. Formatting: header row -> bold font, blue fill, no border, no alignment
$ws.Range("A1:E1").ClearFormats()
$ws.Range("A1:E1").Font.Bold = $true
$ws.Range("A1:E1").Interior.Color = 16711680

# 7. Formatting: data rows 2-7
$ws.Range("A2:A7").ClearFormats()
$ws.Range("A2:A7").Interior.Color = 32768

$ws.Range("B2:B7").ClearFormats()
$ws.Range("B2:B7").Interior.Color = 8421504

$ws.Range("C2:E7").ClearFormats()
$ws.Range("C2:E7").Interior.Color = 65535

$excel.Calculate()
